# The scraper re-run now also pulls "height" and "weight" for each player.
# These two new columns are inserted between the existing "fumbles" column
# (D) and "fantasy points" (which was column E, and is pushed out to G).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank columns at E:F. This shifts the existing "fantasy points"
# column (and all of its data/header formatting) from E to G automatically.
$ws.Range("E1:F1").EntireColumn.Insert()

# New header labels for the inserted columns.
$ws.Range("E1").Value = "height"
$ws.Range("F1").Value = "weight"

# New per-row data: every player row gets the same height/weight values.
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 5).Value = 6.416666666666667
    $ws.Cells.Item($r, 6).Value = 255
}
